# 🔄 Actualización automática del mapa (2025-08-13 08:41:12)
#
# This script reproduces, via Excel COM-interop calls, the edits captured in
# the canonical-XML diff for mapa_interactivo_AYKO.xlsx:
#   1. Row 22 gets its case data corrected (Caso/Direccion/Observaciones/
#      Attachments/coordinates updated).
#   2. A brand-new incident row (Caso 6411, CRAIG 720) is inserted at row 79,
#      pushing the previously-existing rows 79-86 down to 80-87 (dimension
#      grows from A1:P86 to A1:P87).
#   3. One of the shifted rows (now row 85) also needs its "Tipo de
#      Elemento" value corrected from "Pasante" to "Poste".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing it to be stored as TEXT
# (matches the source file, where every column except Attachments/
# Coordenada_X/Coordenada_Y is an inline string -- even when the text
# looks like a plain number, e.g. "6901", "13", "-549", or a date such as
# "2/24/2025"). Without this, Excel's normal type inference would turn
# those values into real numbers / dates.
# ---------------------------------------------------------------------
function Set-TextCell {
    param($Row, $Col, $Text)
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.NumberFormat = "General"
}

function Set-RowData {
    param($Row, [string[]]$Values)
    # Columns: A Caso, B F.De Reclamo, C Direccion, D Comuna, E OT,
    #          F Proveedor Asignado, G Estado, H Observaciones,
    #          I Attachments (numeric), J Tipo de tarea, K Equipo,
    #          L Tipo de Elemento, M Coordenada_X (numeric),
    #          N Coordenada_Y (numeric), O Operacion, P Zona
    Set-TextCell $Row 1  $Values[0]
    Set-TextCell $Row 2  $Values[1]
    Set-TextCell $Row 3  $Values[2]
    Set-TextCell $Row 4  $Values[3]
    Set-TextCell $Row 5  $Values[4]
    Set-TextCell $Row 6  $Values[5]
    Set-TextCell $Row 7  $Values[6]
    Set-TextCell $Row 8  $Values[7]
    $ws.Cells.Item($Row, 9).Value = [double]$Values[8]
    Set-TextCell $Row 10 $Values[9]
    Set-TextCell $Row 11 $Values[10]
    Set-TextCell $Row 12 $Values[11]
    $ws.Cells.Item($Row, 13).Value = [double]$Values[12]
    $ws.Cells.Item($Row, 14).Value = [double]$Values[13]
    Set-TextCell $Row 15 $Values[14]
    Set-TextCell $Row 16 $Values[15]
}

# ---------------------------------------------------------------------
# 1) Update row 22 (case renumbered + address/observation/coord fixes)
# ---------------------------------------------------------------------
Set-TextCell 22 1 "6901"
Set-TextCell 22 3 "PINO, Virrey del 2644"
Set-TextCell 22 8 "Picada"
$ws.Cells.Item(22, 9).Value = 1
$ws.Cells.Item(22, 13).Value = -58.455073
$ws.Cells.Item(22, 14).Value = -34.56776

# ---------------------------------------------------------------------
# 2) Insert the new incident as row 79, shifting the old rows 79-86 down
#    to 80-87 (Excel also grows the used range / dimension to P87).
# ---------------------------------------------------------------------
$ws.Rows(79).Insert()

Set-RowData 79 @(
    "6411",
    "7/31/2025",
    "CRAIG 720",
    "6",
    "808609241",
    "AYKO",
    "Pendiente",
    "Picada",
    "1",
    "Cambio",
    "Sin equipos",
    "Pasante",
    "-58.434949",
    "-34.627171",
    "Boedo",
    "Capital Sur"
)

# ---------------------------------------------------------------------
# 3) The shifted row 85 (old row 84, GOLETA SARANDI 6050) needs its
#    "Tipo de Elemento" corrected from "Pasante" to "Poste".
# ---------------------------------------------------------------------
Set-TextCell 85 12 "Poste"
